$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row right-answer count
$ws.Range("B11").Value = 5

# Update "Total" row right-answer count and the correct/total marks summary
$ws.Range("B12").Value = 105
$ws.Range("E12").Value = "105/140"
